$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122. Used to restore the default (General) style
# after forcing text storage via NumberFormat="@", so numeric-looking
# strings (e.g. "272.94", "-0.04%") are kept as text, matching the
# original inlineStr cells, without leaving a stray cell style behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "272.94"
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.04%"
$ws.Range("E28").Copy()
$ws.Range("E2").PasteSpecial(-4122)

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.72%"
$ws.Range("E28").Copy()
$ws.Range("E3").PasteSpecial(-4122)

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.905"
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.77%"
$ws.Range("E28").Copy()
$ws.Range("E4").PasteSpecial(-4122)

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06321"
$ws.Range("D3").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.25%"
$ws.Range("E28").Copy()
$ws.Range("E5").PasteSpecial(-4122)

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.911"
$ws.Range("D3").Copy()
$ws.Range("D6").PasteSpecial(-4122)

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.15%"
$ws.Range("E28").Copy()
$ws.Range("E6").PasteSpecial(-4122)

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.357"
$ws.Range("D3").Copy()
$ws.Range("D7").PasteSpecial(-4122)

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.28%"
$ws.Range("E28").Copy()
$ws.Range("E7").PasteSpecial(-4122)

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.358"
$ws.Range("D3").Copy()
$ws.Range("D8").PasteSpecial(-4122)

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "49.21%"
$ws.Range("E28").Copy()
$ws.Range("E8").PasteSpecial(-4122)

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8847"
$ws.Range("D3").Copy()
$ws.Range("D9").PasteSpecial(-4122)

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.52%"
$ws.Range("E28").Copy()
$ws.Range("E9").PasteSpecial(-4122)

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.41%"
$ws.Range("E28").Copy()
$ws.Range("E10").PasteSpecial(-4122)

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05120"
$ws.Range("D3").Copy()
$ws.Range("D11").PasteSpecial(-4122)

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.22%"
$ws.Range("E28").Copy()
$ws.Range("E11").PasteSpecial(-4122)

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07393"
$ws.Range("D3").Copy()
$ws.Range("D12").PasteSpecial(-4122)

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.26%"
$ws.Range("E28").Copy()
$ws.Range("E12").PasteSpecial(-4122)

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03118"
$ws.Range("D3").Copy()
$ws.Range("D13").PasteSpecial(-4122)

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.96%"
$ws.Range("E28").Copy()
$ws.Range("E13").PasteSpecial(-4122)

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09037"
$ws.Range("D3").Copy()
$ws.Range("D14").PasteSpecial(-4122)

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.18%"
$ws.Range("E28").Copy()
$ws.Range("E14").PasteSpecial(-4122)

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001557"
$ws.Range("D3").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.80%"
$ws.Range("E28").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006315"
$ws.Range("D3").Copy()
$ws.Range("D16").PasteSpecial(-4122)

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.98%"
$ws.Range("E28").Copy()
$ws.Range("E16").PasteSpecial(-4122)

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006046"
$ws.Range("D3").Copy()
$ws.Range("D17").PasteSpecial(-4122)

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.98%"
$ws.Range("E28").Copy()
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.458"
$ws.Range("D3").Copy()
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.44%"
$ws.Range("E28").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.283"
$ws.Range("D3").Copy()
$ws.Range("D19").PasteSpecial(-4122)

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.79%"
$ws.Range("E28").Copy()
$ws.Range("E19").PasteSpecial(-4122)

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.62%"
$ws.Range("E28").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.63%"
$ws.Range("E28").Copy()
$ws.Range("E21").PasteSpecial(-4122)

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.899"
$ws.Range("D3").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.27%"
$ws.Range("E28").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04369"
$ws.Range("D3").Copy()
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.75%"
$ws.Range("E28").Copy()
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001176"
$ws.Range("D3").Copy()
$ws.Range("D24").PasteSpecial(-4122)

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.10%"
$ws.Range("E28").Copy()
$ws.Range("E24").PasteSpecial(-4122)

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003678"
$ws.Range("D3").Copy()
$ws.Range("D25").PasteSpecial(-4122)

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-12.31%"
$ws.Range("E28").Copy()
$ws.Range("E25").PasteSpecial(-4122)

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("D3").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.08%"
$ws.Range("E28").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "1.34%"
$ws.Range("E28").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04030"
$ws.Range("D3").Copy()
$ws.Range("D40").PasteSpecial(-4122)

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.38%"
$ws.Range("E28").Copy()
$ws.Range("E40").PasteSpecial(-4122)

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006627"
$ws.Range("D3").Copy()
$ws.Range("D41").PasteSpecial(-4122)

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "6.70%"
$ws.Range("E28").Copy()
$ws.Range("E41").PasteSpecial(-4122)

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.85%"
$ws.Range("E28").Copy()
$ws.Range("E42").PasteSpecial(-4122)

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002100"
$ws.Range("D3").Copy()
$ws.Range("D43").PasteSpecial(-4122)

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.14%"
$ws.Range("E28").Copy()
$ws.Range("E43").PasteSpecial(-4122)

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.03%"
$ws.Range("E28").Copy()
$ws.Range("E44").PasteSpecial(-4122)

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005313"
$ws.Range("D3").Copy()
$ws.Range("D45").PasteSpecial(-4122)

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.78%"
$ws.Range("E28").Copy()
$ws.Range("E45").PasteSpecial(-4122)

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.357"
$ws.Range("D3").Copy()
$ws.Range("D46").PasteSpecial(-4122)

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "162.64%"
$ws.Range("E28").Copy()
$ws.Range("E46").PasteSpecial(-4122)

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-12.99%"
$ws.Range("E28").Copy()
$ws.Range("E47").PasteSpecial(-4122)

$excel.CutCopyMode = $false
